# Update the "Estado de Cuenta" worker arrears table (rows 16-41).
# Adds 2 new periods (1810, 1811) for MARCELA MERIÑO OSPINO and re-sorts
# every period block (1812..1905) into the new worker order:
#   YAMADIS CAMARGO MARQUEZ, GINA MARCELA CAMARGO MONROY,
#   OLGA ISABEL FLOREZ MEZA, MARCELA MERIÑO OSPINO
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=16; C="1095811770"; D='MARCELA MERIÑO OSPINO'; E="1810"; F=31249 }
    @{ Row=17; C="1095811770"; D='MARCELA MERIÑO OSPINO'; E="1811"; F=31249 }
    @{ Row=18; C="22815519"; D='YAMADIS CAMARGO MARQUEZ'; E="1812"; F=31249 }
    @{ Row=19; C="1049566193"; D='GINA MARCELA CAMARGO MONROY'; E="1812"; F=31249 }
    @{ Row=20; C="42365405"; D='OLGA ISABEL FLOREZ MEZA'; E="1812"; F=31249 }
    @{ Row=21; C="1095811770"; D='MARCELA MERIÑO OSPINO'; E="1812"; F=31249 }
    @{ Row=22; C="22815519"; D='YAMADIS CAMARGO MARQUEZ'; E="1901"; F=31249 }
    @{ Row=23; C="1049566193"; D='GINA MARCELA CAMARGO MONROY'; E="1901"; F=31249 }
    @{ Row=24; C="42365405"; D='OLGA ISABEL FLOREZ MEZA'; E="1901"; F=31249 }
    @{ Row=25; C="1095811770"; D='MARCELA MERIÑO OSPINO'; E="1901"; F=31249 }
    @{ Row=26; C="22815519"; D='YAMADIS CAMARGO MARQUEZ'; E="1902"; F=31249 }
    @{ Row=27; C="1049566193"; D='GINA MARCELA CAMARGO MONROY'; E="1902"; F=31249 }
    @{ Row=28; C="42365405"; D='OLGA ISABEL FLOREZ MEZA'; E="1902"; F=31249 }
    @{ Row=29; C="1095811770"; D='MARCELA MERIÑO OSPINO'; E="1902"; F=31249 }
    @{ Row=30; C="22815519"; D='YAMADIS CAMARGO MARQUEZ'; E="1903"; F=31249 }
    @{ Row=31; C="1049566193"; D='GINA MARCELA CAMARGO MONROY'; E="1903"; F=31249 }
    @{ Row=32; C="42365405"; D='OLGA ISABEL FLOREZ MEZA'; E="1903"; F=31249 }
    @{ Row=33; C="1095811770"; D='MARCELA MERIÑO OSPINO'; E="1903"; F=31249 }
    @{ Row=34; C="22815519"; D='YAMADIS CAMARGO MARQUEZ'; E="1904"; F=31249 }
    @{ Row=35; C="1049566193"; D='GINA MARCELA CAMARGO MONROY'; E="1904"; F=31249 }
    @{ Row=36; C="42365405"; D='OLGA ISABEL FLOREZ MEZA'; E="1904"; F=31249 }
    @{ Row=37; C="1095811770"; D='MARCELA MERIÑO OSPINO'; E="1904"; F=31249 }
    @{ Row=38; C="22815519"; D='YAMADIS CAMARGO MARQUEZ'; E="1905"; F=20833 }
    @{ Row=39; C="1049566193"; D='GINA MARCELA CAMARGO MONROY'; E="1905"; F=20833 }
    @{ Row=40; C="42365405"; D='OLGA ISABEL FLOREZ MEZA'; E="1905"; F=20833 }
    @{ Row=41; C="1095811770"; D='MARCELA MERIÑO OSPINO'; E="1905"; F=20833 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

